$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country name cells (column A) due to re-ranking / re-sort ---
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 22 de Julio de 2020 a las 01:26'
$ws.Cells.Item(19, 1).Value = 'Colombia'
$ws.Cells.Item(20, 1).Value = 'Banglades'
$ws.Cells.Item(50, 1).Value = 'Nigeria'
$ws.Cells.Item(51, 1).Value = 'Barein'
$ws.Cells.Item(74, 1).Value = 'Venezuela'
$ws.Cells.Item(75, 1).Value = 'El Salvador'
$ws.Cells.Item(76, 1).Value = 'Australia'
$ws.Cells.Item(78, 1).Value = 'Etiopia'
$ws.Cells.Item(79, 1).Value = 'Sudan'
$ws.Cells.Item(104, 1).Value = 'Nicaragua'
$ws.Cells.Item(105, 1).Value = 'Zambia'
$ws.Cells.Item(106, 1).Value = 'Tailandia'
$ws.Cells.Item(127, 1).Value = 'Suazilandia'
$ws.Cells.Item(128, 1).Value = 'Islandia'
$ws.Cells.Item(138, 1).Value = 'Letonia'
$ws.Cells.Item(139, 1).Value = 'Jordania'
$ws.Cells.Item(210, 1).Value = 'Groenlandia'
$ws.Cells.Item(211, 1).Value = 'Islas Malvinas'

# --- Update numeric data cells (columns B-H) with refreshed COVID-19 stats ---
# Row 4
$ws.Cells.Item(4, 2).Value = 4024462
$ws.Cells.Item(4, 3).Value = 63033
$ws.Cells.Item(4, 4).Value = 1884336
$ws.Cells.Item(4, 5).Value = 1995230
$ws.Cells.Item(4, 7).Value = 1062
$ws.Cells.Item(4, 8).Value = 144896
# Row 5
$ws.Cells.Item(5, 2).Value = 2166532
$ws.Cells.Item(5, 3).Value = 44887
$ws.Cells.Item(5, 4).Value = 1465970
$ws.Cells.Item(5, 5).Value = 618965
$ws.Cells.Item(5, 7).Value = 1346
$ws.Cells.Item(5, 8).Value = 81597
# Row 19
$ws.Cells.Item(19, 2).Value = 211038
$ws.Cells.Item(19, 3).Value = 7033
$ws.Cells.Item(19, 4).Value = 98840
$ws.Cells.Item(19, 5).Value = 105032
$ws.Cells.Item(19, 7).Value = 237
$ws.Cells.Item(19, 8).Value = 7166
# Row 20
$ws.Cells.Item(20, 2).Value = 210510
$ws.Cells.Item(20, 3).Value = 3057
$ws.Cells.Item(20, 4).Value = 115397
$ws.Cells.Item(20, 5).Value = 92404
$ws.Cells.Item(20, 7).Value = 41
$ws.Cells.Item(20, 8).Value = 2709
# Row 23
$ws.Cells.Item(23, 2).Value = 136118
$ws.Cells.Item(23, 3).Value = 5344
$ws.Cells.Item(23, 4).Value = 58598
$ws.Cells.Item(23, 5).Value = 75030
$ws.Cells.Item(23, 7).Value = 117
$ws.Cells.Item(23, 8).Value = 2490
# Row 24
$ws.Cells.Item(24, 2).Value = 111669
$ws.Cells.Item(24, 3).Value = 545
$ws.Cells.Item(24, 4).Value = 97742
$ws.Cells.Item(24, 5).Value = 5065
$ws.Cells.Item(24, 7).Value = 4
$ws.Cells.Item(24, 8).Value = 8862
# Row 31
$ws.Cells.Item(31, 2).Value = 76217
$ws.Cells.Item(31, 3).Value = 1597
$ws.Cells.Item(31, 4).Value = 32725
$ws.Cells.Item(31, 5).Value = 38126
$ws.Cells.Item(31, 7).Value = 48
$ws.Cells.Item(31, 8).Value = 5366
# Row 50
$ws.Cells.Item(50, 2).Value = 37801
$ws.Cells.Item(50, 3).Value = 576
$ws.Cells.Item(50, 4).Value = 15677
$ws.Cells.Item(50, 5).Value = 21319
$ws.Cells.Item(50, 7).Value = 4
$ws.Cells.Item(50, 8).Value = 805
# Row 51
$ws.Cells.Item(51, 2).Value = 37316
$ws.Cells.Item(51, 3).Value = 380
$ws.Cells.Item(51, 4).Value = 33455
$ws.Cells.Item(51, 5).Value = 3732
$ws.Cells.Item(51, 7).Value = 1
$ws.Cells.Item(51, 8).Value = 129
# Row 60
$ws.Cells.Item(60, 2).Value = 25736
$ws.Cells.Item(60, 3).Value = 640
$ws.Cells.Item(60, 4).Value = 20155
$ws.Cells.Item(60, 5).Value = 4593
$ws.Cells.Item(60, 7).Value = 3
$ws.Cells.Item(60, 8).Value = 988
# Row 74
$ws.Cells.Item(74, 2).Value = 12774
$ws.Cells.Item(74, 3).Value = 440
$ws.Cells.Item(74, 4).Value = 6983
$ws.Cells.Item(74, 5).Value = 5671
$ws.Cells.Item(74, 7).Value = 4
$ws.Cells.Item(74, 8).Value = 120
# Row 75
$ws.Cells.Item(75, 2).Value = 12582
$ws.Cells.Item(75, 3).Value = 375
$ws.Cells.Item(75, 4).Value = 6965
$ws.Cells.Item(75, 5).Value = 5265
$ws.Cells.Item(75, 7).Value = 8
$ws.Cells.Item(75, 8).Value = 352
# Row 76
$ws.Cells.Item(76, 2).Value = 12428
$ws.Cells.Item(76, 3).Value = 359
$ws.Cells.Item(76, 4).Value = 8444
$ws.Cells.Item(76, 5).Value = 3858
$ws.Cells.Item(76, 7).Value = 3
$ws.Cells.Item(76, 8).Value = 126
# Row 78
$ws.Cells.Item(78, 2).Value = 11072
$ws.Cells.Item(78, 3).Value = 865
$ws.Cells.Item(78, 4).Value = 5448
$ws.Cells.Item(78, 5).Value = 5444
$ws.Cells.Item(78, 7).Value = 10
$ws.Cells.Item(78, 8).Value = 180
# Row 79
$ws.Cells.Item(79, 2).Value = 10992
$ws.Cells.Item(79, 4).Value = 5707
$ws.Cells.Item(79, 5).Value = 4592
$ws.Cells.Item(79, 8).Value = 693
# Row 83
$ws.Cells.Item(83, 2).Value = 9053
$ws.Cells.Item(83, 3).Value = 19
$ws.Cells.Item(83, 5).Value = 660
# Row 93
$ws.Cells.Item(93, 2).Value = 6652
$ws.Cells.Item(93, 3).Value = 62
$ws.Cells.Item(93, 4).Value = 5771
$ws.Cells.Item(93, 5).Value = 840
$ws.Cells.Item(93, 7).Value = 1
$ws.Cells.Item(93, 8).Value = 41
# Row 95
$ws.Cells.Item(95, 2).Value = 5985
$ws.Cells.Item(95, 3).Value = 62
$ws.Cells.Item(95, 4).Value = 3826
$ws.Cells.Item(95, 5).Value = 2004
# Row 98
$ws.Cells.Item(98, 2).Value = 4561
$ws.Cells.Item(98, 3).Value = 13
$ws.Cells.Item(98, 4).Value = 1411
$ws.Cells.Item(98, 5).Value = 3095
# Row 104
$ws.Cells.Item(104, 2).Value = 3439
$ws.Cells.Item(104, 3).Value = 292
$ws.Cells.Item(104, 4).Value = 2492
$ws.Cells.Item(104, 5).Value = 839
$ws.Cells.Item(104, 7).Value = 9
$ws.Cells.Item(104, 8).Value = 108
# Row 105
$ws.Cells.Item(105, 2).Value = 3386
$ws.Cells.Item(105, 3).Value = 60
$ws.Cells.Item(105, 4).Value = 1620
$ws.Cells.Item(105, 5).Value = 1638
$ws.Cells.Item(105, 8).Value = 128
# Row 106
$ws.Cells.Item(106, 2).Value = 3255
$ws.Cells.Item(106, 3).Value = 5
$ws.Cells.Item(106, 4).Value = 3105
$ws.Cells.Item(106, 5).Value = 92
$ws.Cells.Item(106, 8).Value = 58
# Row 127
$ws.Cells.Item(127, 2).Value = 1894
$ws.Cells.Item(127, 3).Value = 68
$ws.Cells.Item(127, 4).Value = 855
$ws.Cells.Item(127, 5).Value = 1015
$ws.Cells.Item(127, 7).Value = 1
$ws.Cells.Item(127, 8).Value = 24
# Row 128
$ws.Cells.Item(128, 2).Value = 1839
$ws.Cells.Item(128, 3).Value = 19
$ws.Cells.Item(128, 4).Value = 1821
$ws.Cells.Item(128, 5).Value = 8
$ws.Cells.Item(128, 7).Value = 0
$ws.Cells.Item(128, 8).Value = 10
# Row 138
$ws.Cells.Item(138, 2).Value = 1193
$ws.Cells.Item(138, 3).Value = 1
$ws.Cells.Item(138, 4).Value = 1045
$ws.Cells.Item(138, 5).Value = 117
$ws.Cells.Item(138, 8).Value = 31
# Row 139
$ws.Cells.Item(139, 2).Value = 1113
$ws.Cells.Item(139, 3).Value = 0
$ws.Cells.Item(139, 4).Value = 1034
$ws.Cells.Item(139, 5).Value = 68
$ws.Cells.Item(139, 8).Value = 11
# Row 181
$ws.Cells.Item(181, 2).Value = 139
$ws.Cells.Item(181, 3).Value = 2
$ws.Cells.Item(181, 5).Value = 7
# Row 182
$ws.Cells.Item(182, 2).Value = 117
$ws.Cells.Item(182, 3).Value = 2
$ws.Cells.Item(182, 5).Value = 14
